$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''63.798.04'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '''  -1.84%  '
$ws.Range('E2').ClearFormats()

$ws.Range('D3').Value = '''3.347.64'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '''  -2.06%  '
$ws.Range('E3').ClearFormats()

$ws.Range('D4').Value = '''1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '''  -0.06%  '
$ws.Range('E4').ClearFormats()

$ws.Range('D5').Value = '''546.19'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '''  +0.15%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').Value = '''172.29'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '''  -3.23%  '
$ws.Range('E6').ClearFormats()

$ws.Range('E7').Value = '''  -3.05%  '
$ws.Range('E7').ClearFormats()

$ws.Range('D8').Value = '''3.334.91'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '''  -2.17%  '
$ws.Range('E8').ClearFormats()

$ws.Range('D9').Value = '''1.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '''  +0.04%  '
$ws.Range('E9').ClearFormats()

$ws.Range('D10').Value = '''0.613'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '''  -0.96%  '
$ws.Range('E10').ClearFormats()

$ws.Range('E11').Value = '''  +1.98%  '
$ws.Range('E11').ClearFormats()

$ws.Range('D12').Value = '''53.93'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '''  +1.85%  '
$ws.Range('E12').ClearFormats()

$ws.Range('D13').Value = '''0.0000265'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '''  -1.06%  '
$ws.Range('E13').ClearFormats()

$ws.Range('D14').Value = '''8.91'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '''  -1.97%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').Value = '''3.881.74'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '''  -2.50%  '
$ws.Range('E15').ClearFormats()

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '''3.370.23'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '''  -1.69%  '
$ws.Range('E16').ClearFormats()

$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '''17.93'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '''  -1.05%  '
$ws.Range('E17').ClearFormats()

$ws.Range('E18').Value = '''  -3.00%  '
$ws.Range('E18').ClearFormats()

$ws.Range('D19').Value = '''11.71'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '''  +0.03%  '
$ws.Range('E19').ClearFormats()

$ws.Range('D20').Value = '''63.816.40'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '''  -2.10%  '
$ws.Range('E20').ClearFormats()

$ws.Range('D21').Value = '''0.976'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '''  +0.27%  '
$ws.Range('E21').ClearFormats()

$ws.Range('D22').Value = '''411.91'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '''  -0.01%  '
$ws.Range('E22').ClearFormats()

$ws.Range('E23').Value = '''  +1.72%  '
$ws.Range('E23').ClearFormats()

$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '''4.33'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '''  +1.77%  '
$ws.Range('E24').ClearFormats()

$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = '''13.76'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '''  +14.35%  '
$ws.Range('E25').ClearFormats()

$ws.Range('D26').Value = '''83.02'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '''  -1.37%  '
$ws.Range('E26').ClearFormats()

$ws.Range('D27').Value = '''10.56'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '''  -1.47%  '
$ws.Range('E27').ClearFormats()

$ws.Range('E28').Value = '''  -3.53%  '
$ws.Range('E28').ClearFormats()

$ws.Range('D29').Value = '''8.60'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '''  -2.32%  '
$ws.Range('E29').ClearFormats()

$ws.Range('D30').Value = '''29.05'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '''  -1.59%  '
$ws.Range('E30').ClearFormats()

$ws.Range('D31').Value = '''6.40'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '''  -0.50%  '
$ws.Range('E31').ClearFormats()

$ws.Range('D32').Value = '''11.33'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '''  -1.78%  '
$ws.Range('E32').ClearFormats()

$ws.Range('D33').Value = '''569.87'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '''  -6.55%  '
$ws.Range('E33').ClearFormats()

$ws.Range('E34').Value = '''  -1.24%  '
$ws.Range('E34').ClearFormats()

$ws.Range('D35').Value = '''57.93'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '''  -0.96%  '
$ws.Range('E35').ClearFormats()

$ws.Range('E36').Value = '''  +0.60%  '
$ws.Range('E36').ClearFormats()

$ws.Range('D37').Value = '''0.998'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '''  +0.08%  '
$ws.Range('E37').ClearFormats()

$ws.Range('D38').Value = '''35.07'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '''  -5.02%  '
$ws.Range('E38').ClearFormats()

$ws.Range('D39').Value = '''3.39'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '''  +2.40%  '
$ws.Range('E39').ClearFormats()

$ws.Range('D40').Value = '''0.0₃0738'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '''  -4.05%  '
$ws.Range('E40').ClearFormats()

$ws.Range('D41').Value = '''0.367'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '''  -1.90%  '
$ws.Range('E41').ClearFormats()

$ws.Range('D42').Value = '''3.150.85'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '''  +0.16%  '
$ws.Range('E42').ClearFormats()

$ws.Range('D43').Value = '''1.00'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '''  -0.24%  '
$ws.Range('E43').ClearFormats()

$ws.Range('E44').Value = '''  +1.43%  '
$ws.Range('E44').ClearFormats()

$ws.Range('D45').Value = '''3.26'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '''  +2.71%  '
$ws.Range('E45').ClearFormats()

$ws.Range('E46').Value = '''  -1.12%  '
$ws.Range('E46').ClearFormats()

$ws.Range('D47').Value = '''2.41'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '''  -4.36%  '
$ws.Range('E47').ClearFormats()

$ws.Range('E48').Value = '''  -4.49%  '
$ws.Range('E48').ClearFormats()

$ws.Range('E49').Value = '''  -1.64%  '
$ws.Range('E49').ClearFormats()

$ws.Range('D50').Value = '''132.34'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '''  -4.08%  '
$ws.Range('E50').ClearFormats()

$ws.Range('D51').Value = '''8.07'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '''  -2.83%  '
$ws.Range('E51').ClearFormats()

